# Updated library selection script: the SamplesTab query's row-limit was
# reduced from 100 to 10 (global / limit fix), and the "SamplesTab" query
# cell (B3) is now the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$oldText = $ws.Range("B3").Value2
$newText = $oldText.Replace("ORDER By samp.sample_id LIMIT 100", "ORDER By samp.sample_id LIMIT 10")
$ws.Range("B3").Value = $newText

$ws.Range("B3").Select()
